# Applies updated detailed-accuracy values (re-run results after adding a
# dropdown to select the source video) to Sheet1 of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = 99.99985694885254
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 99.99984502792358
$ws.Range("D5").Value = 100
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 98.8919198513031
$ws.Range("D7").Value = 99.99996423721313
$ws.Range("D8").Value = 99.99983310699463
$ws.Range("D9").Value = 99.42722320556641
$ws.Range("D10").Value = 99.99996423721313
$ws.Range("D16").Value = 99.99986886978149
$ws.Range("D17").Value = 99.99974966049194
$ws.Range("D18").Value = 99.99995231628418
$ws.Range("D19").Value = 100
$ws.Range("D20").Value = 97.98550009727478
$ws.Range("D22").Value = 99.99998807907104
$ws.Range("D23").Value = 99.89905953407288
$ws.Range("D26").Value = 99.99731779098511
$ws.Range("D27").Value = 99.99998807907104
$ws.Range("D29").Value = 99.99998807907104
$ws.Range("D32").Value = 100
$ws.Range("D33").Value = 99.99998807907104
$ws.Range("D34").Value = 99.99990463256836
$ws.Range("D35").Value = 99.99830722808838
$ws.Range("D39").Value = 99.99927282333374
$ws.Range("D40").Value = 99.99998807907104
$ws.Range("D41").Value = 99.99788999557495
$ws.Range("D43").Value = 99.78626370429993
$ws.Range("D44").Value = 91.97914004325867
$ws.Range("D45").Value = 99.99995231628418
$ws.Range("D46").Value = 99.99994039535522
$ws.Range("C48").Value = 1
$ws.Range("D48").Value = 91.61857962608337
$ws.Range("D50").Value = 99.97105002403259
$ws.Range("D52").Value = 99.99972581863403
$ws.Range("D53").Value = 100
$ws.Range("D54").Value = 99.99943971633911
$ws.Range("D55").Value = 100
$ws.Range("D59").Value = 99.988853931427
$ws.Range("D61").Value = 99.96284246444702
$ws.Range("C63").Value = 2
$ws.Range("D63").Value = 99.99997615814209
$ws.Range("D64").Value = 99.99961853027344
